$d = $word.ActiveDocument

# Remove the "SeatNo." property line from the FlightTicket data design section.
# The text currently reads "...FlightClass<br/>SeatNo.<br/>..." and must become
# "...FlightClass<br/>..." (the SeatNo line, including its leading line break,
# is deleted).
$r = $d.Content
$found = $r.Find.Execute("SeatNo." + [char]11, $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
if ($found) {
    $r.Delete()
}
